$d = $word.ActiveDocument

$replacements = @(
    @{old = "76×96="; new = "29×55="},
    @{old = "61×25="; new = "17×69="},
    @{old = "21×17="; new = "16×91="},
    @{old = "36×99="; new = "63×23="},
    @{old = "53×67="; new = "31×92="},
    @{old = "49×65="; new = "57×30="},
    @{old = "32×86="; new = "72×36="},
    @{old = "21×32="; new = "46×84="},
    @{old = "35×39="; new = "24×84="},
    @{old = "80×67="; new = "37×49="},
    @{old = "40×72="; new = "99×21="},
    @{old = "81×30="; new = "25×58="},
    @{old = "86×98="; new = "46×28="},
    @{old = "34×66="; new = "21×17="},
    @{old = "14×71="; new = "76×39="},
    @{old = "83×93="; new = "99×83="},
    @{old = "25×71="; new = "65×89="},
    @{old = "76×90="; new = "50×79="},
    @{old = "21×23="; new = "89×92="},
    @{old = "85×85="; new = "87×17="},
    @{old = "38×59="; new = "77×32="},
    @{old = "24×59="; new = "80×98="},
    @{old = "82×67="; new = "29×93="},
    @{old = "65×46="; new = "85×47="},
    @{old = "57×21="; new = "53×64="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
